$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 44, shifting the existing rows 44:140 down to 45:141
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new data entry (new weekly record
# for "Ají" - Americana (o), added ahead of the previously-first record)
$ws.Cells.Item(44, 1).Value = 5
$ws.Cells.Item(44, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(44, 3).Value = 'Maule'
$ws.Cells.Item(44, 4).Value = 44536
$ws.Cells.Item(44, 5).Value = 7
$ws.Cells.Item(44, 6).Value = 100112021
$ws.Cells.Item(44, 7).Value = 'Ají'
$ws.Cells.Item(44, 8).Value = 'Americana (o)'
$ws.Cells.Item(44, 9).Value = 'Primera'
$ws.Cells.Item(44, 10).Value = 150
$ws.Cells.Item(44, 11).Value = 18000
$ws.Cells.Item(44, 12).Value = 18000
$ws.Cells.Item(44, 13).Value = 18000
$ws.Cells.Item(44, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(44, 15).Value = 'Región del Maule'
$ws.Cells.Item(44, 16).Value = 1200
$ws.Cells.Item(44, 17).Value = 15
$ws.Cells.Item(44, 18).Value = 'Hortaliza'
